# Update forests data - 2025-10-02 12:16
#
# Sheet "New" currently holds 3 freshly scraped listings (rows 2-4).
# This edit:
#   1) Moves those 3 existing listings (with their hyperlinks) down into
#      the "Previously added" sheet, appended as new rows 144-146.
#   2) Replaces the "New" sheet's rows 2-4 with 3 brand-new listings
#      (values + hyperlinks).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Previously added")
$ws2 = $wb.Worksheets.Item("New")

# Helper: write literal text into a cell without Excel's "looks like a
# number/date" auto-coercion (and without leaving a quote-prefix style
# behind). We do this by writing a text formula and then collapsing it
# to a literal value in place.
function Set-TextValue($range, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
    $excel.CutCopyMode = $false
}

# ---------------------------------------------------------------------
# Step 1: append the current "New" rows 2-4 onto "Previously added" as
# rows 144-146, carrying over formatting, values (shared strings reused
# as-is) and hyperlinks.
# ---------------------------------------------------------------------

# Clone formatting from the last existing data row (143) down onto the
# three new rows, so styles (s="3"/"4"/"2") match the rest of the table.
$ws1.Range("A143:F143").Copy()
$ws1.Range("A144:F146").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Copy the actual cell values over (text cells keep their shared-string
# type, numeric/date cell F keeps its numeric type).
$ws2.Range("A2:F4").Copy()
$ws1.Range("A144:F146").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

# Re-create the hyperlinks on column A for the 3 appended rows, pointing
# at the same URLs as the source cells on the "New" sheet.
$ws1.Hyperlinks.Add($ws1.Range("A144"), $ws1.Range("A144").Value2)
$ws1.Hyperlinks.Add($ws1.Range("A145"), $ws1.Range("A145").Value2)
$ws1.Hyperlinks.Add($ws1.Range("A146"), $ws1.Range("A146").Value2)

# Hyperlinks.Add() auto-applies the built-in "Hyperlink" cell style;
# re-assert the workbook's own link style (matching every other A-column
# cell in this sheet) on top of it.
$ws1.Range("A143").Copy()
$ws1.Range("A144").PasteSpecial(-4122)  # xlPasteFormats
$ws1.Range("A145").PasteSpecial(-4122)
$ws1.Range("A146").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Step 2: overwrite "New" rows 2-4 with the 3 newly scraped listings.
# ---------------------------------------------------------------------

# Drop the old hyperlinks on A2:A4 before writing new ones on top.
$ws2.Range("A2:A4").Hyperlinks.Delete()

# Row 2
Set-TextValue $ws2.Range("A2") "https://www.ss.com/msg/lv/real-estate/wood/kraslava-and-reg/udrisu-pag/bkbix.html"
Set-TextValue $ws2.Range("B2") "35 880 €"
Set-TextValue $ws2.Range("C2") "Krāslava un raj."
Set-TextValue $ws2.Range("D2") "8 ha."
Set-TextValue $ws2.Range("E2") "60960020067/68"
$ws2.Range("F2").Value = 45932.393055555556

# Row 3
Set-TextValue $ws2.Range("A3") "https://www.ss.com/msg/lv/real-estate/wood/limbadzi-and-reg/salacgrivas-l-t/cdkbp.html"
Set-TextValue $ws2.Range("B3") "50 000 €"
Set-TextValue $ws2.Range("C3") "Limbaži un raj."
Set-TextValue $ws2.Range("D3") "3 ha."
Set-TextValue $ws2.Range("E3") "66720040252"
$ws2.Range("F3").Value = 45931.83125

# Row 4
Set-TextValue $ws2.Range("A4") "https://www.ss.com/msg/lv/real-estate/wood/limbadzi-and-reg/liepupes-pag/eepmh.html"
Set-TextValue $ws2.Range("B4") "39 000 €"
Set-TextValue $ws2.Range("C4") "Limbaži un raj."
Set-TextValue $ws2.Range("D4") "6.50 ha."
Set-TextValue $ws2.Range("E4") "66600090044"
$ws2.Range("F4").Value = 45931.757638888885

# Re-add hyperlinks for the new rows
$ws2.Hyperlinks.Add($ws2.Range("A2"), $ws2.Range("A2").Value2)
$ws2.Hyperlinks.Add($ws2.Range("A3"), $ws2.Range("A3").Value2)
$ws2.Hyperlinks.Add($ws2.Range("A4"), $ws2.Range("A4").Value2)

# Hyperlinks.Add() auto-applies the built-in "Hyperlink" cell style;
# re-assert the sheet's own link style (style index 3, same as before
# the overwrite) on the A-column cells.
$ws1.Range("A143").Copy()
$ws2.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$ws2.Range("A3").PasteSpecial(-4122)
$ws2.Range("A4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

Write-Output "forests data updated"
